{"js": "// Remove 4 consecutive empty \"ACL Text First Line\" paragraphs that\n// immediately follow the \"Tabelle training\" paragraph (the first one,\n// exactly matching that text, not \"Tabelle training Anf\u00e4nger\" etc.).\n// There were originally 14 such empty paragraphs in a row; after the\n// edit only 10 remain.\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"text,style\");\nawait context.sync();\n\n// Locate the \"Tabelle training\" paragraph that is followed by a run of\n// empty \"ACL Text First Line\" paragraphs (disambiguates from the other\n// \"Tabelle training ...\" headings elsewhere in the document).\nlet targetIndex = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"Tabelle training\") {\n    let emptyCount = 0;\n    let j = i + 1;\n    while (\n      j < paras.items.length &&\n      paras.items[j].text === \"\" &&\n      paras.items[j].style === \"ACL Text First Line\"\n    ) {\n      emptyCount++;\n      j++;\n    }\n    if (emptyCount >= 4) {\n      targetIndex = i;\n      break;\n    }\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error('Could not find \"Tabelle training\" paragraph followed by empty lines.');\n}\n\n// Delete the first 4 empty paragraphs right after the target paragraph.\nfor (let k = 0; k < 4; k++) {\n  paras.items[targetIndex + 1 + k].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove 4 consecutive empty \"ACL Text First Line\" paragraphs that\n# immediately follow the \"Tabelle training\" paragraph (the first one,\n# exactly matching that text, not \"Tabelle training Anf\u00e4nger\" etc.).\n# There were originally 14 such empty paragraphs in a row; after the\n# edit only 10 remain.\n\n$d = $word.ActiveDocument\n$wdParagraph = 4\n$styleName = \"ACL Text First Line\"\n\n$count = $d.Paragraphs.Count\n\nfunction Test-EmptyStyledParagraph($p, $name) {\n  if ($p -eq $null) { return $false }\n  $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n  return ($t -eq \"\" -and $p.Style.NameLocal -eq $name)\n}\n\n# Walk every \"Tabelle training\" occurrence (Find) until we find the one\n# that is immediately followed by >= 4 empty \"ACL Text First Line\"\n# paragraphs - this disambiguates it from the other \"Tabelle training ...\"\n# table captions elsewhere in the document.\n$searchFrom = 0\n$docEnd = $d.Content.End\n$targetIdx = -1\n\nwhile ($true) {\n  $rng = $d.Range($searchFrom, $docEnd)\n  $rng.Find.ClearFormatting()\n  $rng.Find.Text = \"Tabelle training\"\n  $rng.Find.Forward = $true\n  $rng.Find.Wrap = 0\n  $ok = $rng.Find.Execute()\n  if (-not $ok) { break }\n\n  $para = $rng.Duplicate\n  $null = $para.Expand($wdParagraph)\n  $paraText = $para.Text.TrimEnd([char]13, [char]7)\n\n  if ($paraText -eq \"Tabelle training\") {\n    # Resolve this paragraph's index in $d.Paragraphs by matching\n    # Start/End/Text (not just Start - some paragraphs can report\n    # identical cached Start offsets).\n    $idx = -1\n    for ($i = 1; $i -le $count; $i++) {\n      $p = $d.Paragraphs.Item($i)\n      if ($p.Range.Start -eq $para.Start -and $p.Range.End -eq $para.End -and $p.Range.Text -eq $para.Text) {\n        $idx = $i\n        break\n      }\n    }\n\n    if ($idx -ne -1) {\n      $p1 = $d.Paragraphs.Item($idx + 1)\n      $p2 = $d.Paragraphs.Item($idx + 2)\n      $p3 = $d.Paragraphs.Item($idx + 3)\n      $p4 = $d.Paragraphs.Item($idx + 4)\n\n      if ((Test-EmptyStyledParagraph $p1 $styleName) -and\n          (Test-EmptyStyledParagraph $p2 $styleName) -and\n          (Test-EmptyStyledParagraph $p3 $styleName) -and\n          (Test-EmptyStyledParagraph $p4 $styleName)) {\n        $targetIdx = $idx\n        break\n      }\n    }\n  }\n\n  $searchFrom = $rng.End\n  if ($searchFrom -ge $docEnd) { break }\n}\n\nif ($targetIdx -eq -1) {\n  throw \"Could not locate the 'Tabelle training' paragraph followed by empty lines.\"\n}\n\n# Delete the first 4 empty paragraphs right after the target paragraph.\n$p1 = $d.Paragraphs.Item($targetIdx + 1)\n$p4 = $d.Paragraphs.Item($targetIdx + 4)\n$delRange = $d.Range($p1.Range.Start, $p4.Range.End)\n$delRange.Delete()\n"}
